$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old header row outright - this drops its bold/centered/wrapText
# style and the taller row height along with it (rows 2-4 shift up to 1-3).
$ws.Rows.Item(1).Delete()

# Drop the now-unused extra data row (old row 4, now row 3) so only two data
# rows remain, then insert a fresh row 1 for the new header.
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(1).Insert()

# Rebuild the table as grade_id / school_id / grade_name
$ws.Range("A1").Value = "grade_id"
$ws.Range("B1").Value = "school_id"
$ws.Range("C1").Value = "grade_name"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "Grade 1"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Grade 2"

# Column B no longer needs its custom width; column C gets one instead.
$ws.Columns.Item(2).ColumnWidth = 7.59
$ws.Columns.Item(3).ColumnWidth = 10.59

# Match the recorded cursor position after the edit.
$ws.Range("A11").Select()
